{"js": "// Apply the Spanish translation fixes described by the diff:\n// 1. \"[ID de cuenta]\" -> \"[account ID]\"\n// 2. \"...antes de retirar su saldo.\" -> \"...antes de retirar su balance.\" (2nd occurrence of \"su saldo\")\n// 3. \"se cerrar\u00e1 el 29/09/2023\" -> \"se cerrar\u00e1 \u00e9l 29/09/2023\"\n// 4. \"el saldo de la cuenta se transferir\u00e1\" -> \"el balance de la cuenta se transferir\u00e1\"\n// 5. \"Durante este proceso se aplicar\u00e1n los tipos de cambio y las tasas habituales.\"\n//    -> \"Durante este proceso se aplicar\u00e1n las tasas de cambio y las tarifas habituales.\"\n\nconst body = context.document.body;\n\nconst replacements = [\n  { find: \"[ID de cuenta]\", replace: \"[account ID]\" },\n  {\n    find: \"ci\u00e9rrelas primero antes de retirar su saldo.\",\n    replace: \"ci\u00e9rrelas primero antes de retirar su balance.\",\n  },\n  {\n    find: \"Su cuenta USDT se cerrar\u00e1 el 29/09/2023\",\n    replace: \"Su cuenta USDT se cerrar\u00e1 \u00e9l 29/09/2023\",\n  },\n  {\n    find: \"y el saldo de la cuenta se transferir\u00e1\",\n    replace: \"y el balance de la cuenta se transferir\u00e1\",\n  },\n  {\n    find: \"Durante este proceso se aplicar\u00e1n los tipos de cambio y las tasas habituales.\",\n    replace: \"Durante este proceso se aplicar\u00e1n las tasas de cambio y las tarifas habituales.\",\n  },\n];\n\nfor (const { find, replace } of replacements) {\n  const results = body.search(find, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${find}`);\n  }\n\n  for (const r of results.items) {\n    r.insertText(replace, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the Spanish translation fixes described by the diff:\n# 1. \"[ID de cuenta]\" -> \"[account ID]\"\n# 2. \"...antes de retirar su saldo.\" -> \"...antes de retirar su balance.\"\n# 3. \"se cerrar\u00e1 el 29/09/2023\" -> \"se cerrar\u00e1 \u00e9l 29/09/2023\"\n# 4. \"el saldo de la cuenta se transferir\u00e1\" -> \"el balance de la cuenta se transferir\u00e1\"\n# 5. \"Durante este proceso se aplicar\u00e1n los tipos de cambio y las tasas habituales.\"\n#    -> \"Durante este proceso se aplicar\u00e1n las tasas de cambio y las tarifas habituales.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 0          # wdFindStop\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$findText, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]0, [ref]$false, [ref]$replaceText, [ref]2) | Out-Null\n}\n\nReplace-Text \"[ID de cuenta]\" \"[account ID]\"\nReplace-Text \"ci\u00e9rrelas primero antes de retirar su saldo.\" \"ci\u00e9rrelas primero antes de retirar su balance.\"\nReplace-Text \"Su cuenta USDT se cerrar\u00e1 el 29/09/2023\" \"Su cuenta USDT se cerrar\u00e1 \u00e9l 29/09/2023\"\nReplace-Text \"y el saldo de la cuenta se transferir\u00e1\" \"y el balance de la cuenta se transferir\u00e1\"\nReplace-Text \"Durante este proceso se aplicar\u00e1n los tipos de cambio y las tasas habituales.\" \"Durante este proceso se aplicar\u00e1n las tasas de cambio y las tarifas habituales.\"\n"}
